$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update testpass -> password in B2 and B3
$ws.Range("B2").Value = "password"
$ws.Range("B3").Value = "password"

# Update selection to reflect new active cell
$ws.Range("M11").Select()
